$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Col4a1","Itga1","ECs",3,1,216.5832213333333,649.7496640000001,0.4331411212367192,0.4331411212367192,3,1,43.28121633333333,129.843649,0.533749049291363,0.533749049291363,9373.985256698215,84365.86731028394,0.2311886616690939,0.2311886616690939),
  @("ECs","Col4a1","Itga1","FAPs",3,1,216.5832213333333,649.7496640000001,0.4331411212367192,0.4331411212367192,3,1,9.865038,29.595114,0.1216568086527629,0.1216568086527629,2136.601708615744,19229.4153775417,0.05269456650593873,0.05269456650593871),
  @("ECs","Col4a1","Itga1","M2",3,1,216.5832213333333,649.7496640000001,0.4331411212367192,0.4331411212367192,2,0.6666666666666666,1.460025,4.380075,0.01800519998536753,0.01800519998536753,316.2169177272,2845.9522595448,0.007798792509753453,0.007798792509753452),
  @("ECs","Col4a1","Itga1","sCs",3,1,216.5832213333333,649.7496640000001,0.4331411212367192,0.4331411212367192,3,1,26.482795,79.448385,0.3265889420705065,0.3265889420705065,5735.729051010293,51621.56145909265,0.1414591005519331,0.1414591005519331),
  @("FAPs","Col4a1","Itga1","ECs",3,1,161.954974,485.864922,0.3238910133313607,0.3238910133313606,3,1,43.28121633333333,129.843649,0.533749049291363,0.533749049291363,7009.608265953375,63086.47439358038,0.1728765204396299,0.1728765204396299),
  @("FAPs","Col4a1","Itga1","FAPs",3,1,161.954974,485.864922,0.3238910133313607,0.3238910133313606,3,1,9.865038,29.595114,0.1216568086527629,0.1216568086527629,1597.691972799012,14379.22775519111,0.03940354703320283,0.03940354703320281),
  @("FAPs","Col4a1","Itga1","M2",3,1,161.954974,485.864922,0.3238910133313607,0.3238910133313606,2,0.6666666666666666,1.460025,4.380075,0.01800519998536753,0.01800519998536753,236.45831091435,2128.12479822915,0.00583172246849449,0.005831722468494489),
  @("FAPs","Col4a1","Itga1","sCs",3,1,161.954974,485.864922,0.3238910133313607,0.3238910133313606,3,1,26.482795,79.448385,0.3265889420705065,0.3265889420705065,4289.020375672329,38601.18338105097,0.1057792233900334,0.1057792233900334),
  @("M2","Col4a1","Itga1","ECs",3,1,0.4608033333333333,1.38241,0.0009215527926904059,0.0009215527926904059,3,1,43.28121633333333,129.843649,0.533749049291363,0.533749049291363,19.94412875712111,179.49715881409,0.0004918779269703047,0.0004918779269703047),
  @("M2","Col4a1","Itga1","FAPs",3,1,0.4608033333333333,1.38241,0.0009215527926904059,0.0009215527926904059,3,1,9.865038,29.595114,0.1216568086527629,0.1216568086527629,4.545842393859999,40.91258154473999,0.000112113171763756,0.000112113171763756),
  @("M2","Col4a1","Itga1","M2",3,1,0.4608033333333333,1.38241,0.0009215527926904059,0.0009215527926904059,2,0.6666666666666666,1.460025,4.380075,0.01800519998536753,0.01800519998536753,0.6727843867499999,6.055059480749999,0.0000165927423294647,0.0000165927423294647),
  @("M2","Col4a1","Itga1","sCs",3,1,0.4608033333333333,1.38241,0.0009215527926904059,0.0009215527926904059,3,1,26.482795,79.448385,0.3265889420705065,0.3265889420705065,12.20336021198333,109.83024190785,0.0003009689516268805,0.0003009689516268805),
  @("sCs","Col4a1","Itga1","ECs",3,1,121.0302313333333,363.090694,0.2420463126392298,0.2420463126392298,3,1,43.28121633333333,129.843649,0.533749049291363,0.533749049291363,5238.335625211378,47145.02062690241,0.1291919892556689,0.1291919892556689),
  @("sCs","Col4a1","Itga1","FAPs",3,1,121.0302313333333,363.090694,0.2420463126392298,0.2420463126392298,3,1,9.865038,29.595114,0.1216568086527629,0.1216568086527629,1193.967831252124,10745.71048126912,0.02944658194185761,0.0294465819418576),
  @("sCs","Col4a1","Itga1","M2",3,1,121.0302313333333,363.090694,0.2420463126392298,0.2420463126392298,2,0.6666666666666666,1.460025,4.380075,0.01800519998536753,0.01800519998536753,176.70716350245,1590.36447152205,0.004358092264790125,0.004358092264790125),
  @("sCs","Col4a1","Itga1","sCs",3,1,121.0302313333333,363.090694,0.2420463126392298,0.2420463126392298,3,1,26.482795,79.448385,0.3265889420705065,0.3265889420705065,3205.218805203243,28846.96924682919,0.07904964917691314,0.07904964917691312)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowArr = $data[$r]
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value = $rowArr[$c]
    }
}

